$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: WP4754 / IL-18 signaling pathway (metrics refreshed) ---
$ws.Cells.Item(2, 1).Value = "WP4754"
$ws.Cells.Item(2, 2).Value = "IL-18 signaling pathway"
$ws.Cells.Item(2, 3).Value = 2.0
$ws.Cells.Item(2, 4).Value = -0.828125
$ws.Cells.Item(2, 5).Value = -1.365903289497467
$ws.Cells.Item(2, 6).Value = 0.0759493670886076
$ws.Cells.Item(2, 7).Value = 0.45569620253164556
$ws.Cells.Item(2, 8).Value = 0.45569620253164556
$ws.Cells.Item(2, 9).Value = 14.0
$ws.Cells.Item(2, 10).Value = "tags=100%, list=21%, signal=81%"
$ws.Cells.Item(2, 11).Value = "'7078"

# --- Row 3: WP2879 / Farnesoid X Receptor  Pathway ---
$ws.Cells.Item(3, 1).Value = "WP2879"
$ws.Cells.Item(3, 2).Value = "Farnesoid X Receptor  Pathway"
$ws.Cells.Item(3, 3).Value = 2.0
$ws.Cells.Item(3, 4).Value = 0.78125
$ws.Cells.Item(3, 5).Value = 1.3571337353949977
$ws.Cells.Item(3, 6).Value = 0.1206140350877193
$ws.Cells.Item(3, 7).Value = 0.45569620253164556
$ws.Cells.Item(3, 8).Value = 0.45569620253164556
$ws.Cells.Item(3, 9).Value = 16.0
$ws.Cells.Item(3, 10).Value = "tags=100%, list=24%, signal=78%"
$ws.Cells.Item(3, 11).Value = "2289/117283"

# --- Row 4: WP2877 / Vitamin D Receptor Pathway ---
$ws.Cells.Item(4, 1).Value = "WP2877"
$ws.Cells.Item(4, 2).Value = "Vitamin D Receptor Pathway"
$ws.Cells.Item(4, 3).Value = 2.0
$ws.Cells.Item(4, 4).Value = -0.796875
$ws.Cells.Item(4, 5).Value = -1.3143597691390718
$ws.Cells.Item(4, 6).Value = 0.1301989150090416
$ws.Cells.Item(4, 7).Value = 0.45569620253164556
$ws.Cells.Item(4, 8).Value = 0.45569620253164556
$ws.Cells.Item(4, 9).Value = 16.0
$ws.Cells.Item(4, 10).Value = "tags=100%, list=24%, signal=78%"
$ws.Cells.Item(4, 11).Value = "'7078"

# --- Row 5: WP1449 / Regulation of toll-like receptor signaling pathway ---
$ws.Cells.Item(5, 1).Value = "WP1449"
$ws.Cells.Item(5, 2).Value = "Regulation of toll-like receptor signaling pathway"
$ws.Cells.Item(5, 3).Value = 2.0
$ws.Cells.Item(5, 4).Value = -0.7781510728319369
$ws.Cells.Item(5, 5).Value = -1.2834766612614346
$ws.Cells.Item(5, 6).Value = 0.1518987341772152
$ws.Cells.Item(5, 7).Value = 0.45569620253164556
$ws.Cells.Item(5, 8).Value = 0.45569620253164556
$ws.Cells.Item(5, 9).Value = 2.0
$ws.Cells.Item(5, 10).Value = "tags=100%, list=3%, signal=100%"
$ws.Cells.Item(5, 11).Value = "'23098"

# --- Row 6 (WP3932 / Focal Adhesion-PI3K-Akt-mTOR-signaling pathway) is removed ---
$ws.Rows(6).Delete()
